# Business Exception logging update:
# Append a new logged exception (row 7) and four more rows (8-11) that
# repeat earlier exception messages, mirroring the shared-strings/rows
# added in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Exceptions")

$newText = "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-22-2020\RandomFolder2\ for the change: CHANGES - SOX Audit Report for magic435.txt_07.01.73.eml made on 1/22/2020 is not a valid path."

# Row 7 previously held no data; it now carries the newly logged exception
# text (this becomes shared-string index 3).
$ws.Range("A7").Value2 = $newText

# Rows 8-11 repeat earlier messages (shared-string indices 0, 2, 0, 1).
$ws.Range("A8").Value2  = $ws.Range("A1").Value2
$ws.Range("A9").Value2  = $ws.Range("A5").Value2
$ws.Range("A10").Value2 = $ws.Range("A1").Value2
$ws.Range("A11").Value2 = $ws.Range("A2").Value2

# Carry over the existing cell formatting (style index 1) onto the newly
# populated rows, matching the rest of column A.
$ws.Range("A6").Copy()
$ws.Range("A7:A11").PasteSpecial(-4122)
